# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Hyperion_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 32.153847
$ws.Range("I5").Value = 32.153847
$ws.Range("K5").Value = 32.153847
$ws.Range("M5").Value = 82.846153
$ws.Range("H70").Value = 1443.6666
$ws.Range("I70").Value = 999
$ws.Range("J70").Value = 1499.25
$ws.Range("K70").Value = 2997
$ws.Range("L70").Value = 4497.75
$ws.Range("M70").Value = -2727
$ws.Range("N70").Value = -5037.75
$ws.Range("H73").Value = 1443.6666
$ws.Range("I73").Value = 999
$ws.Range("J73").Value = 1499.25
$ws.Range("K73").Value = 2997
$ws.Range("L73").Value = 4497.75
$ws.Range("M73").Value = -2061
$ws.Range("N73").Value = -6369.75
$ws.Range("H98").Value = 6728.625
$ws.Range("I98").Value = 6728.625
$ws.Range("K98").Value = 6728.625
$ws.Range("M98").Value = -5230.625
$ws.Range("H115").Value = 330.88235
$ws.Range("I115").Value = 330.88235
$ws.Range("K115").Value = 992.6470499999999
$ws.Range("M115").Value = 574.3529500000001
$ws.Range("H122").Value = 6728.625
$ws.Range("I122").Value = 6728.625
$ws.Range("K122").Value = 20185.875
$ws.Range("M122").Value = -17735.875
$ws.Range("H132").Value = 13892993
$ws.Range("I132").Value = 17547826
$ws.Range("J132").Value = 4624.533
$ws.Range("K132").Value = 52643478
$ws.Range("L132").Value = 13873.599
$ws.Range("M132").Value = -52640948
$ws.Range("N132").Value = -18933.599
$ws.Range("H136").Value = 135935.33
$ws.Range("J136").Value = 135935.33
$ws.Range("L136").Value = 135935.33
$ws.Range("N136").Value = -146135.33
$ws.Range("H137").Value = 58440.195
$ws.Range("I137").Value = 75027.03999999999
$ws.Range("J137").Value = 1571
$ws.Range("K137").Value = 225081.12
$ws.Range("L137").Value = 4713
$ws.Range("M137").Value = -222531.12
$ws.Range("N137").Value = -9813
$ws.Range("H138").Value = 4574.9165
$ws.Range("I138").Value = 3225
$ws.Range("J138").Value = 5094.115
$ws.Range("K138").Value = 9675
$ws.Range("L138").Value = 15282.345
$ws.Range("M138").Value = -4535
$ws.Range("N138").Value = -25562.345

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 298.75
$ws.Range("I4").Value = 298.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 298.75
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -182.75
$ws.Range("N4").Value = ""
$ws.Range("H5").Value = 5832.5
$ws.Range("I5").Value = 328.2857
$ws.Range("K5").Value = 328.2857
$ws.Range("M5").Value = -216.2857
$ws.Range("H25").Value = 1513.625
$ws.Range("I25").Value = 1515.5714
$ws.Range("J25").Value = 1500
$ws.Range("K25").Value = 1515.5714
$ws.Range("L25").Value = 1500
$ws.Range("M25").Value = -1113.5714
$ws.Range("N25").Value = -2304
$ws.Range("H32").Value = 3303.9868
$ws.Range("I32").Value = 2218.5386
$ws.Range("K32").Value = 2218.5386
$ws.Range("M32").Value = -1931.5386
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = ""
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -54992
$ws.Range("H97").Value = 2029452
$ws.Range("I97").Value = 2489874.8
$ws.Range("K97").Value = 2489874.8
$ws.Range("M97").Value = -2489378.8
$ws.Range("H122").Value = 455561.56
$ws.Range("I122").Value = 2314.0908
$ws.Range("K122").Value = 6942.2724
$ws.Range("M122").Value = -4492.2724
$ws.Range("H132").Value = 3288.9727
$ws.Range("I132").Value = 3100.9167
$ws.Range("K132").Value = 9302.750100000001
$ws.Range("M132").Value = -6772.750100000001
$ws.Range("H138").Value = 51750
$ws.Range("J138").Value = 51750
$ws.Range("L138").Value = 51750
$ws.Range("N138").Value = -62030

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 35000
$ws.Range("J2").Value = 35000
$ws.Range("L2").Value = 35000
$ws.Range("N2").Value = -35226
$ws.Range("H4").Value = 5832.5
$ws.Range("I4").Value = 328.2857
$ws.Range("K4").Value = 328.2857
$ws.Range("M4").Value = -213.2857
$ws.Range("H22").Value = 824.8333
$ws.Range("I22").Value = 718
$ws.Range("K22").Value = 718
$ws.Range("M22").Value = -545
$ws.Range("H86").Value = 4174656
$ws.Range("I86").Value = 5892476
$ws.Range("J86").Value = 2807.7144
$ws.Range("K86").Value = 5892476
$ws.Range("L86").Value = 2807.7144
$ws.Range("M86").Value = -5891353
$ws.Range("N86").Value = -5053.7144
$ws.Range("H89").Value = 4174656
$ws.Range("I89").Value = 5892476
$ws.Range("J89").Value = 2807.7144
$ws.Range("K89").Value = 29462380
$ws.Range("L89").Value = 14038.572
$ws.Range("M89").Value = -29456764
$ws.Range("N89").Value = -25270.572
$ws.Range("H92").Value = 45100.25
$ws.Range("J92").Value = 45100.25
$ws.Range("L92").Value = 45100.25
$ws.Range("N92").Value = -50092.25
$ws.Range("H94").Value = 10104639
$ws.Range("I94").Value = 30307252
$ws.Range("J94").Value = 3333
$ws.Range("K94").Value = 30307252
$ws.Range("L94").Value = 3333
$ws.Range("M94").Value = -30306801
$ws.Range("N94").Value = -4235
$ws.Range("H110").Value = 14000
$ws.Range("J110").Value = 14000
$ws.Range("L110").Value = 14000
$ws.Range("N110").Value = -22180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2749.6
$ws.Range("I16").Value = 2666.2222
$ws.Range("K16").Value = 2666.2222
$ws.Range("M16").Value = -2379.2222
$ws.Range("H31").Value = 3610.2122
$ws.Range("I31").Value = 1168.2858
$ws.Range("J31").Value = 4267.654
$ws.Range("K31").Value = 1168.2858
$ws.Range("L31").Value = 4267.654
$ws.Range("M31").Value = -873.2858000000001
$ws.Range("N31").Value = -4857.654
$ws.Range("H34").Value = 3610.2122
$ws.Range("I34").Value = 1168.2858
$ws.Range("J34").Value = 4267.654
$ws.Range("K34").Value = 1168.2858
$ws.Range("L34").Value = 4267.654
$ws.Range("M34").Value = -966.2858000000001
$ws.Range("N34").Value = -4671.654
$ws.Range("H58").Value = 2271.1428
$ws.Range("I58").Value = 2615.2856
$ws.Range("K58").Value = 2615.2856
$ws.Range("M58").Value = -2412.2856
$ws.Range("H62").Value = 2998.75
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 2998.3333
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 2998.3333
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4246.3333
$ws.Range("H65").Value = 2998.75
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 2998.3333
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 14991.6665
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -21231.6665
$ws.Range("H93").Value = 15985
$ws.Range("J93").Value = 50874.5
$ws.Range("L93").Value = 50874.5
$ws.Range("N93").Value = -54618.5
$ws.Range("H113").Value = 2749.6
$ws.Range("I113").Value = 2666.2222
$ws.Range("K113").Value = 2666.2222
$ws.Range("M113").Value = -496.2222000000002
$ws.Range("H132").Value = 43773.918
$ws.Range("I132").Value = 51954
$ws.Range("J132").Value = 2873.5
$ws.Range("K132").Value = 155862
$ws.Range("L132").Value = 8620.5
$ws.Range("M132").Value = -153332
$ws.Range("N132").Value = -13680.5
$ws.Range("H134").Value = 9263.532999999999
$ws.Range("I134").Value = 6335.522
$ws.Range("K134").Value = 19006.566
$ws.Range("M134").Value = -16471.566
$ws.Range("H136").Value = 2271.1428
$ws.Range("I136").Value = 2615.2856
$ws.Range("K136").Value = 7845.8568
$ws.Range("M136").Value = -5295.8568
$ws.Range("H138").Value = 93937.5
$ws.Range("J138").Value = 93937.5
$ws.Range("L138").Value = 93937.5
$ws.Range("N138").Value = -104217.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 65363.93
$ws.Range("J12").Value = 2904.2222
$ws.Range("L12").Value = 8712.6666
$ws.Range("N12").Value = -9058.6666
$ws.Range("H68").Value = 1522.8108
$ws.Range("I68").Value = 1299.8334
$ws.Range("J68").Value = 1734.0526
$ws.Range("K68").Value = 3899.5002
$ws.Range("L68").Value = 5202.1578
$ws.Range("M68").Value = -3088.5002
$ws.Range("N68").Value = -6824.1578
$ws.Range("H71").Value = 1522.8108
$ws.Range("I71").Value = 1299.8334
$ws.Range("J71").Value = 1734.0526
$ws.Range("K71").Value = 11698.5006
$ws.Range("L71").Value = 15606.4734
$ws.Range("M71").Value = -7642.500599999999
$ws.Range("N71").Value = -23718.4734
$ws.Range("H121").Value = 627.9091
$ws.Range("I121").Value = 272.4
$ws.Range("J121").Value = 924.1667
$ws.Range("K121").Value = 817.1999999999999
$ws.Range("L121").Value = 2772.5001
$ws.Range("M121").Value = 492.8000000000001
$ws.Range("N121").Value = -5392.5001
$ws.Range("H136").Value = 2525.6
$ws.Range("I136").Value = 2282
$ws.Range("K136").Value = 6846
$ws.Range("M136").Value = -1746
$ws.Range("H137").Value = 3402.76
$ws.Range("I137").Value = 2459.7
$ws.Range("J137").Value = 4031.4666
$ws.Range("K137").Value = 7379.099999999999
$ws.Range("L137").Value = 12094.3998
$ws.Range("M137").Value = -2279.099999999999
$ws.Range("N137").Value = -22294.3998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6331834
$ws.Range("I80").Value = 11319309
$ws.Range("K80").Value = 11319309
$ws.Range("M80").Value = -11318311
$ws.Range("H83").Value = 6331834
$ws.Range("I83").Value = 11319309
$ws.Range("K83").Value = 56596545
$ws.Range("M83").Value = -56591553
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H102").Value = 3091403.8
$ws.Range("I102").Value = 3833262.5
$ws.Range("K102").Value = 3833262.5
$ws.Range("M102").Value = -3831640.5
$ws.Range("H122").Value = 333390.03
$ws.Range("I122").Value = 448731.6
$ws.Range("K122").Value = 1346194.8
$ws.Range("M122").Value = -1343744.8
$ws.Range("H126").Value = 2262071.8
$ws.Range("I126").Value = 1198541.9
$ws.Range("J126").Value = 4389131
$ws.Range("K126").Value = 3595625.7
$ws.Range("L126").Value = 13167393
$ws.Range("M126").Value = -3593155.7
$ws.Range("N126").Value = -13172333
$ws.Range("H132").Value = 4799
$ws.Range("I132").Value = 2555.6487
$ws.Range("K132").Value = 7666.946100000001
$ws.Range("M132").Value = -5136.946100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3399.1482
$ws.Range("I7").Value = 2064
$ws.Range("K7").Value = 2064
$ws.Range("M7").Value = -1952
$ws.Range("H17").Value = 2547.75
$ws.Range("I17").Value = 1730.3334
$ws.Range("K17").Value = 1730.3334
$ws.Range("M17").Value = -1560.3334
$ws.Range("H19").Value = 579.1667
$ws.Range("I19").Value = 751.5
$ws.Range("J19").Value = 493
$ws.Range("K19").Value = 751.5
$ws.Range("L19").Value = 493
$ws.Range("M19").Value = -581.5
$ws.Range("N19").Value = -833
$ws.Range("H22").Value = 51177.89
$ws.Range("I22").Value = 149430
$ws.Range("J22").Value = 2051.8333
$ws.Range("K22").Value = 149430
$ws.Range("L22").Value = 2051.8333
$ws.Range("M22").Value = -149135
$ws.Range("N22").Value = -2641.8333
$ws.Range("H27").Value = 51177.89
$ws.Range("I27").Value = 149430
$ws.Range("J27").Value = 2051.8333
$ws.Range("K27").Value = 149430
$ws.Range("L27").Value = 2051.8333
$ws.Range("M27").Value = -149323
$ws.Range("N27").Value = -2265.8333
$ws.Range("H40").Value = 2638091
$ws.Range("I40").Value = 4353243
$ws.Range("J40").Value = 8191.7334
$ws.Range("K40").Value = 4353243
$ws.Range("L40").Value = 8191.7334
$ws.Range("M40").Value = -4353107
$ws.Range("N40").Value = -8463.733400000001
$ws.Range("H46").Value = 3771.2856
$ws.Range("I46").Value = 1900
$ws.Range("J46").Value = 4519.8
$ws.Range("K46").Value = 1900
$ws.Range("L46").Value = 4519.8
$ws.Range("M46").Value = -1712
$ws.Range("N46").Value = -4895.8
$ws.Range("H61").Value = 15875366
$ws.Range("I61").Value = 23811624
$ws.Range("K61").Value = 23811624
$ws.Range("M61").Value = -23811422
$ws.Range("H68").Value = 3222
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 3222
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H93").Value = 33337028
$ws.Range("I93").Value = 41669670
$ws.Range("J93").Value = 6450
$ws.Range("K93").Value = 41669670
$ws.Range("L93").Value = 6450
$ws.Range("M93").Value = -41668422
$ws.Range("N93").Value = -8946
$ws.Range("H100").Value = 5512.375
$ws.Range("J100").Value = 5728.4287
$ws.Range("L100").Value = 5728.4287
$ws.Range("N100").Value = -6810.4287
$ws.Range("H113").Value = 15875366
$ws.Range("I113").Value = 23811624
$ws.Range("K113").Value = 23811624
$ws.Range("M113").Value = -23809454
$ws.Range("H122").Value = 5359.826
$ws.Range("I122").Value = 3605.9285
$ws.Range("K122").Value = 10817.7855
$ws.Range("M122").Value = -8367.7855
$ws.Range("H126").Value = 3399.1482
$ws.Range("I126").Value = 2064
$ws.Range("K126").Value = 6192
$ws.Range("M126").Value = -3722
$ws.Range("H132").Value = 2916.1667
$ws.Range("I132").Value = 2741.1562
$ws.Range("K132").Value = 8223.4686
$ws.Range("M132").Value = -5693.4686

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 228.33333
$ws.Range("I14").Value = 217.5
$ws.Range("J14").Value = 250
$ws.Range("K14").Value = 217.5
$ws.Range("L14").Value = 250
$ws.Range("M14").Value = -49.5
$ws.Range("N14").Value = -586
$ws.Range("H45").Value = 10522.857
$ws.Range("J45").Value = 10522.857
$ws.Range("L45").Value = 10522.857
$ws.Range("N45").Value = -11504.857
$ws.Range("H55").Value = 3991.6667
$ws.Range("I55").Value = 5500
$ws.Range("K55").Value = 5500
$ws.Range("M55").Value = -5223
$ws.Range("H94").Value = 24299.4
$ws.Range("J94").Value = 24299.4
$ws.Range("L94").Value = 24299.4
$ws.Range("N94").Value = -26101.4
$ws.Range("H96").Value = 4067
$ws.Range("I96").Value = 3880.4
$ws.Range("K96").Value = 3880.4
$ws.Range("M96").Value = -2507.4
$ws.Range("H113").Value = 1108.4517
$ws.Range("I113").Value = 1132.8572
$ws.Range("K113").Value = 3398.5716
$ws.Range("M113").Value = -1228.5716
$ws.Range("H122").Value = 2034.6316
$ws.Range("I122").Value = 1953.8928
$ws.Range("K122").Value = 5861.678400000001
$ws.Range("M122").Value = -3411.678400000001
$ws.Range("H126").Value = 3958.9656
$ws.Range("I126").Value = 3646.739
$ws.Range("K126").Value = 10940.217
$ws.Range("M126").Value = -8470.217000000001
$ws.Range("H132").Value = 25928242
$ws.Range("I132").Value = 32264270
$ws.Range("J132").Value = 1376135.4
$ws.Range("K132").Value = 96792810
$ws.Range("L132").Value = 4128406.2
$ws.Range("M132").Value = -96790280
$ws.Range("N132").Value = -4133466.2
$ws.Range("H136").Value = 7837.5713
$ws.Range("I136").Value = 10489.272
$ws.Range("J136").Value = 3350.077
$ws.Range("K136").Value = 31467.816
$ws.Range("L136").Value = 10050.231
$ws.Range("M136").Value = -28917.816
$ws.Range("N136").Value = -15150.231
$ws.Range("H141").Value = 124497.5
$ws.Range("J141").Value = 124497.5
$ws.Range("L141").Value = 124497.5
$ws.Range("N141").Value = -134857.5
